$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right above the current row 670 (Betarraga / Primera+Segunda
# weekly entries). This pushes the existing rows 670-740 down to 672-742, which keeps
# every other row's data intact (Excel shifts cell contents automatically).
$ws.Rows("670:671").Insert()

# Row 670: new "Primera" quality entry
$ws.Cells.Item(670, 1).Value = 9
$ws.Cells.Item(670, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(670, 3).Value = "Metropolitana"
$ws.Cells.Item(670, 4).Value = 44946
$ws.Cells.Item(670, 5).Value = 13
$ws.Cells.Item(670, 6).Value = 100114014
$ws.Cells.Item(670, 7).Value = "Betarraga"
$ws.Cells.Item(670, 8).Value = "Sin especificar"
$ws.Cells.Item(670, 9).Value = "Primera"
$ws.Cells.Item(670, 10).Value = 10600
$ws.Cells.Item(670, 11).Value = 80
$ws.Cells.Item(670, 12).Value = 90
$ws.Cells.Item(670, 13).Value = 85
$ws.Cells.Item(670, 14).Value = "$/unidad"
$ws.Cells.Item(670, 15).Value = "Región Metropolitana"
$ws.Cells.Item(670, 16).Value = 85
$ws.Cells.Item(670, 17).Value = 1
$ws.Cells.Item(670, 18).Value = "Hortaliza"

# Row 671: new "Segunda" quality entry
$ws.Cells.Item(671, 1).Value = 9
$ws.Cells.Item(671, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(671, 3).Value = "Metropolitana"
$ws.Cells.Item(671, 4).Value = 44946
$ws.Cells.Item(671, 5).Value = 13
$ws.Cells.Item(671, 6).Value = 100114014
$ws.Cells.Item(671, 7).Value = "Betarraga"
$ws.Cells.Item(671, 8).Value = "Sin especificar"
$ws.Cells.Item(671, 9).Value = "Segunda"
$ws.Cells.Item(671, 10).Value = 5200
$ws.Cells.Item(671, 11).Value = 70
$ws.Cells.Item(671, 12).Value = 70
$ws.Cells.Item(671, 13).Value = 70
$ws.Cells.Item(671, 14).Value = "$/unidad"
$ws.Cells.Item(671, 15).Value = "Región Metropolitana"
$ws.Cells.Item(671, 16).Value = 70
$ws.Cells.Item(671, 17).Value = 1
$ws.Cells.Item(671, 18).Value = "Hortaliza"
